# "Añadido esqueleto y primera prueba para el controlador del tablero"
#
# Adds a new "Tablero" test section (class header / method / single test
# case) at the bottom of the existing test-case table, mirroring the
# "Bot" / "Red Neuronal Controller" sections already present, and widens
# column C so the longer scenario text fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C is widened to fit the new, longer scenario description ---
$ws.Columns("C").ColumnWidth = 40.25

# --- New section header row: "Tablero" (same look as "Bot" / "Red Neuronal Controller") ---
$ws.Range("A2:F2").Copy()
$ws.Range("A56:F56").PasteSpecial(-4122)
$ws.Range("A56").Value = "Tablero"

# --- Method name row: "Crear tablero" ---
$ws.Range("B57").Value = "Crear tablero"
$ws.Range("E4").Copy()
$ws.Range("E57").PasteSpecial(-4122)

# --- First test case for the new method ---
$ws.Range("C58").Value = "Todo vacío excepto las manos de los jugadores"
$ws.Range("C58").WrapText = $false
$ws.Range("E4").Copy()
$ws.Range("E58").PasteSpecial(-4122)
$ws.Range("E58").Formula = "=E55+1"
$ws.Range("F5").Copy()
$ws.Range("F58").PasteSpecial(-4122)
$ws.Range("F58").Value = "Correcto"

# --- The previously-selected leftover formatted cell D63 moves out of the
#     way of the new rows; D64 becomes the new trailing placeholder cell ---
$ws.Range("D63").Clear()
$ws.Range("D64").Font.Name = "Calibri"

# --- F45 loses the stray underline formatting it had ---
$ws.Range("F45").Font.Underline = -4142

# --- Final selection left on D64, matching the saved workbook state ---
$ws.Range("D64").Select()
